# Weekly refresh of the "Hortaliza, Terminal La Palmera de La Serena - Espinaca"
# sheet: a new week's record is prepended at row 18 (the data rows are sorted
# with the newest entries pushed toward the top of the block starting at row
# 18), shifting every existing row from 18..194 down by one (to 19..195).
#
# All columns except D (Fecha), J (Volumen), K (Precio minimo), M (Precio
# promedio ponderado) and P (Precio $/Kg) are constant across every data row
# in this sheet, so the inserted row reuses those constant values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 18:194 down to 19:195, carrying formatting (incl. the date
# style on column D) along with them.
$ws.Rows("18:18").Insert()

# Populate the newly inserted row 18 with the new week's data.
$ws.Range("A18").Value = 8
$ws.Range("B18").Value = "Terminal La Palmera de La Serena"
$ws.Range("C18").Value = "Coquimbo"
$ws.Range("D18").Value = 44545
$ws.Range("E18").Value = 4
$ws.Range("F18").Value = 100112012
$ws.Range("G18").Value = "Espinaca"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 2800
$ws.Range("K18").Value = 450
$ws.Range("L18").Value = 500
$ws.Range("M18").Value = 475
$ws.Range("N18").Value = "`$/atado 300 a 500 gramos"
$ws.Range("O18").Value = "Provincia del Elquí"
$ws.Range("P18").Value = 950
$ws.Range("Q18").Value = 0.5
$ws.Range("R18").Value = "Hortaliza"
